# Updates the benchmark stats table to match the refreshed README.md
# figures for the Renaissance / JDK17 / ShenandoahGC "scala-doku" run.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Single-value cells (rows 1-12).
$t.Cell(1, 1).Range.Text  = "0M"
$t.Cell(2, 1).Range.Text  = "0M"
$t.Cell(3, 1).Range.Text  = "0M"
$t.Cell(4, 1).Range.Text  = "1404"
$t.Cell(5, 1).Range.Text  = "0.00001"
$t.Cell(6, 1).Range.Text  = "0.00289"
$t.Cell(7, 1).Range.Text  = "0.00018"
$t.Cell(8, 1).Range.Text  = "0.00007"
$t.Cell(9, 1).Range.Text  = "0.00029"
$t.Cell(10, 1).Range.Text = "0.00037"
$t.Cell(11, 1).Range.Text = "0.00042"
$t.Cell(12, 1).Range.Text = "0.29479"

# Rows that previously packed an entire tab-delimited summary line into a
# single run now collapse down to just the first (headline) figure.
$t.Cell(44, 1).Range.Text = "99.9"
$t.Cell(45, 1).Range.Text = "0.29"
$t.Cell(46, 1).Range.Text = "293"
